# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the data source.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 2343
    $ws.Range("F3").Value = 1829
    $ws.Range("F5").Value = 1124
    $ws.Range("F6").Value = 1054
    $ws.Range("F8").Value = 5927
}

$wb.Save()
